$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Path" header in A1, matching the style already used by B1 (bold, bordered, centered)
$ws.Range("A1").Value = "Path"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update IPC execution time values in column B (rows 2-31)
$values = @(
    2.487842556667488,
    2.49074680761321,
    2.387410944463444,
    2.381252895741406,
    2.077716119191529,
    2.073477671924134,
    2.089003324738202,
    2.088002545352321,
    2.118193458516476,
    2.125606919169309,
    2.16028060440714,
    2.139607640507212,
    2.365426027882202,
    2.294032299898257,
    2.300343855353062,
    2.078483567790478,
    2.1169956413529,
    2.068225692092291,
    2.066854092554542,
    2.09851492523725,
    2.105226826843285,
    2.279505060492302,
    2.295387219711106,
    2.276925417356452,
    2.227327930608984,
    2.082177426834825,
    2.073989146208124,
    2.066287352165516,
    1.777900528508751,
    1.744290951808829
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
